# The document has a "提交时间：" (submission date) line that currently reads
# "提交时间：           2019年 4月 17日" and an "实习单位（盖章）：" (unit seal) line
# whose 33 trailing underlined spaces are split into two runs (2 + 31) with the
# "_GoBack" bookmark sitting between them.
#
# The edit blanks out the date digits (2019 / 4 / 17), replacing them with
# plain spaces of the same visual width, and relocates the "_GoBack" bookmark
# from the unit-seal line to sit right after the (now blank) year field on the
# date line - which also collapses the unit-seal line's two space runs back
# into a single 33-space run once the bookmark is no longer splitting them.

$d = $word.ActiveDocument

# --- 1. Blank the "2019" year value and drop the relocated bookmark right after it ---
$yearRange = $d.Content
$found = $yearRange.Find.Execute("2019", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $yearRange.Text = "  "
    $bmRange = $d.Range($yearRange.End, $yearRange.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- 2. Blank the "4" month value (the digit immediately before "月") ---
$monthRange = $d.Content
$found = $monthRange.Find.Execute("4月", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $monthDigit = $d.Range($monthRange.Start, $monthRange.Start + 1)
    $monthDigit.Text = "   "
}

# --- 3. Blank the " 17" day value (leading space plus the two digits) ---
$dayRange = $d.Content
$found = $dayRange.Find.Execute(" 17", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $dayRange.Text = "    "
}

# --- 4. Re-merge the unit-seal line's trailing spaces into a single run ---
$labelRange = $d.Content
$found = $labelRange.Find.Execute("实习单位（盖章）：", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $labelEnd = $labelRange.End

    # Locate the paragraph that holds the label, without assuming its index.
    $targetPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $labelEnd -and $labelEnd -le $p.Range.End) {
            $targetPara = $p
            break
        }
    }

    if ($targetPara -ne $null) {
        $paraEnd = $targetPara.Range.End - 1   # exclude the paragraph mark
        $tailRange = $d.Range($labelEnd, $paraEnd)

        # Force a genuine content change first so the engine actually rebuilds
        # the run (writing back the identical text is treated as a no-op and
        # leaves the original two runs - split by the old bookmark - in place).
        $tailRange.Text = "X"
        $paraEndAfter = $targetPara.Range.End - 1
        $tailRange2 = $d.Range($labelEnd, $paraEndAfter)
        $tailRange2.Text = "                                 "
    }
}

Write-Output "edit complete"
